$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sales data update: quantidade (quantity) sold for each row is now 1
# (reflects the updated vendas.csv / vendas.xlsx source data for FIFA FC 25).
$ws.Range("D2:D10").Value = 1

# Scroll the view down a bit and move the selection, matching the
# author's on-screen state when the file was saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1

$ws.Range("D11").Select()
